$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row heights shrink back to 75 for a few existing wrapped-text rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(23).RowHeight = 75
$ws.Rows.Item(27).RowHeight = 75
$ws.Rows.Item(62).RowHeight = 75

# ---------------------------------------------------------------------------
# 2) Row 100 (id 10126, preferredLang / Arabic description) had its
#    value_json corrected: the placeholder English word "arabic" is replaced
#    by the actual Arabic translation, while the lang_code stays "ara".
# ---------------------------------------------------------------------------
$ws.Range("E100").Value = '{"value":"عربي","code":"ara"}'

# ---------------------------------------------------------------------------
# 3) Three new preferredLang rows are appended (101-103). Each is seeded by
#    copying row 99's cell formatting (it already carries the wrap-text
#    style on B/C/E and the boolean style on G that the new rows need), then
#    the data_type cell (D) - which the new rows leave blank - is cleared
#    and the real values are written in.
# ---------------------------------------------------------------------------

# --- Row 101: English description, value_json now in Arabic, lang eng ----
$ws.Range("A99:I99").Copy($ws.Range("A101:I101"))
$ws.Range("D101").ClearContents()
$ws.Range("F101").WrapText = $true
$ws.Rows.Item(101).RowHeight = 30

$ws.Range("A101").Value = 10127
$ws.Range("B101").Value = "preferredLang"
$ws.Range("C101").Value = "user preferred Language"
$ws.Range("E101").Value = '{"value":"عربي","code":"ara"}'
$ws.Range("F101").Value = "eng"

# --- Row 102: French description, value_json in Arabic, lang fra ---------
$ws.Range("A99:I99").Copy($ws.Range("A102:I102"))
$ws.Range("D102").ClearContents()
$ws.Range("F102").WrapText = $true
$ws.Rows.Item(102).RowHeight = 30

$ws.Range("A102").Value = 10128
$ws.Range("B102").Value = "preferredLang"
$ws.Range("C102").Value = "Langue préférée de l'utilisateur"
$ws.Range("E102").Value = '{"value":"عربي","code":"ara"}'
$ws.Range("F102").Value = "fra"

# --- Row 103: Arabic description, value_json in French, lang ara ---------
$ws.Range("A99:I99").Copy($ws.Range("A103:I103"))
$ws.Range("D103").ClearContents()
$ws.Range("F103").WrapText = $true

$ws.Range("A103").Value = 10129
$ws.Range("B103").Value = "preferredLang"
$ws.Range("C103").Value = "يفضل المستخدم اللغة"
$ws.Range("E103").Value = '{"value":"français","code":"fra"}'
$ws.Range("F103").Value = "ara"

# ---------------------------------------------------------------------------
# 4) AutoFilter now covers the grown table (A1:I103 instead of A1:I67).
#    Toggling Range.AutoFilter() on the currently-active range switches it
#    off first, then re-enabling it with the new range turns it back on.
# ---------------------------------------------------------------------------
$ws.Range("A1:I67").AutoFilter()
$ws.Range("A1:I103").AutoFilter()

# The hidden _FilterDatabase defined name needs to track the same range.
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Sheet1!`$A`$1:`$I`$103"

# ---------------------------------------------------------------------------
# 5) Selection / scroll position ends up parked near the newly-added rows.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
$ws.Range("F101").Select()
